$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Thünen-Institut entry (row 36): short German/English names gain
# the "Johann Heinrich von" prefix.
$ws.Range("D36").Value = "Johann Heinrich von Thünen-Institut"
$ws.Range("E36").Value = "Johann Heinrich von Thünen Institute"

# Add a new source row 53 (Umweltbundesamt / German Environment Agency).
# First copy formatting from the row above (52) so the new row matches the
# existing table styling (borders/fill/font per column).
$ws.Range("A52:J52").Copy()
$ws.Range("A53:J53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A53").Value = "Q_UBA_1"
$ws.Range("B53").Value = "Umweltbundesamt nach Angaben der Länder und Flussgebietsgemeinschaften"
$ws.Range("C53").Value = "German Environment Agency (as reported by the Länder and by river basin commissions"
$ws.Range("D53").Value = "Umweltbundesamt nach Angaben der Länder und Flussgebietsgemeinschaften"
$ws.Range("E53").Value = "German Environment Agency (as reported by the Länder and by river basin commissions"
$ws.Range("F53").Value = "https://www.umweltbundesamt.de/"
$ws.Range("G53").Value = "https://www.umweltbundesamt.de/en"
$ws.Range("H53").Value = ""
$ws.Range("I53").Value = ""
$ws.Range("J53").Value = "uba"
